$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-22 down to 12-23
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new weekly data entry
$ws.Cells.Item(11, 1).Value = 11
$ws.Cells.Item(11, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(11, 3).Value = "Bíobío"
$ws.Cells.Item(11, 4).Value = 44533
$ws.Cells.Item(11, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = 100112026
$ws.Cells.Item(11, 7).Value = "Haba"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 180
$ws.Cells.Item(11, 11).Value = 8000
$ws.Cells.Item(11, 12).Value = 8500
$ws.Cells.Item(11, 13).Value = 8222
$ws.Cells.Item(11, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(11, 15).Value = "Región del Maule"
$ws.Cells.Item(11, 16).Value = 329
$ws.Cells.Item(11, 17).Value = 25
$ws.Cells.Item(11, 18).Value = "Hortaliza"
